$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: rename the open/close pin header pair for the two-hall layout ---
$ws.Range("I3").Value = "O2"
$ws.Range("J3").Value = "O1"
$ws.Range("K3").Value = "C1"
$ws.Range("L3").Value = "C2"

# --- Row 7: two new labelled cells describing close-out wiring per hall ---
$ws.Range("N7").Value = "C2O"
$ws.Range("O7").Value = "C1O"

# --- Legend block: close-in pins (rows 11-12) ---
$ws.Range("A11").Value = "CI1"
$ws.Range("B11").Value = "Close In First hall"
$ws.Range("A12").Value = "CI2"
$ws.Range("B12").Value = "Close In Second hall"

# --- Legend block: open-in pins (rows 13-14) ---
$ws.Range("A13").Value = "OI1"
$ws.Range("B13").Value = "Open in First Hall"
$ws.Range("A14").Value = "OI2"
$ws.Range("B14").Value = "Open in Second hall"

# --- Legend block: close-out pins (rows 16-17) ---
$ws.Range("A16").Value = "CO1"
$ws.Range("B16").Value = "Close Out First hall"
$ws.Range("A17").Value = "CO2"
$ws.Range("B17").Value = "Close Out Second hall"

# --- Old row 18 (BU / Button Up) moves down, and the button legend gets its
#     own three rows (19-21) below the blank separator row 18 ---
$ws.Range("A18:B18").ClearContents()

$ws.Range("A19").Value = "BD"
$ws.Range("B19").Value = "Button Down"
$ws.Range("A20").Value = "BM"
$ws.Range("B20").Value = "Button Middle"
$ws.Range("A21").Value = "BU"
$ws.Range("B21").Value = "Button Up"

# --- Selection moves to the new O1/O2/C1/C2 header range ---
$ws.Range("I2:L2").Select()
